# Adds two new columns, I ("I0") and J ("IF"), to the existing data table.
# Header cells (row 1) get the same bold/bordered style used by the other
# header cells (copied from H1); data cells (rows 2-24) stay unstyled like
# the rest of the numeric columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) -------------------------------------------------
# Copy H1's formatting (bold font + border) onto the two new header cells
# before writing their text so they match the rest of the header row.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data (rows 2-24) --------------------------------------------------
$data = @(
    @(2, 7, 7),
    @(3, 7, 7),
    @(4, 8, 8),
    @(5, 6, 7),
    @(6, 5, 5),
    @(7, 7, 7),
    @(8, 8, 8),
    @(9, 6, 6),
    @(10, 8, 8),
    @(11, 9, 9),
    @(12, 4, 5),
    @(13, 8, 8),
    @(14, 7, 7),
    @(15, 10, 10),
    @(16, 7, 7),
    @(17, 7, 8),
    @(18, 7, 8),
    @(19, 7, 8),
    @(20, 6, 6),
    @(21, 5, 5),
    @(22, 8, 8),
    @(23, 7, 7),
    @(24, 3, 3)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]

    $ws.Cells.Item($row, 9).Value = $iVal   # column I
    $ws.Cells.Item($row, 10).Value = $jVal  # column J
}
